$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "3100580400000V"
$ws.Range("B2").Value = "SOCIETE RWANDAISE DES PNEUMATIQUES BANDAG  LTD"
$ws.Range("D2").Value = "SOCIETE RWANDAISE DES PNEUMATIQUES BANDAG  LTD"
$ws.Range("E2").Value = "Kicukiro- Rwanda"
$ws.Range("F2").Value = "Kicukiro- Rwanda"
$ws.Range("G2").Value = "Kicukiro- Rwanda"
$ws.Range("H2").Value = "1132 kigali"
$ws.Range("J2").Value = "bandag@subizo.com"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "0788303361"
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = "AA3035975"
$ws.Range("O2").Value = "ITALY"
$ws.Range("P2").Value = "I&M Bank"
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "100003368"

# Row 4
$ws.Range("A4").Value = "3100960800000N"
$ws.Range("B4").Value = "INYANGE INDUSTRIES  LTD"
$ws.Range("D4").Value = "INYANGE INDUSTRIES  LTD"
$ws.Range("E4").Value = "MASAKA GASABO DISTRICT"
$ws.Range("F4").Value = "MASAKA GASABO DISTRICT"
$ws.Range("G4").Value = "MASAKA GASABO DISTRICT"
$ws.Range("H4").Value = "4584 kigali-rwanda"
$ws.Range("J4").Value = "bjames@inyangeindustries.com"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "0788309662"
$ws.Range("M4").Value = "www.inyangeindustries.com"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = "1197580006310064"
$ws.Range("O4").Value = "Rwandan"
$ws.Range("P4").Value = "NCBA"
$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "100095380"

# Row 5
$ws.Range("A5").Value = "3102226900000M"
$ws.Range("B5").Value = "BROADBAND SYSTEMS CORPORATION  LTD"
$ws.Range("D5").Value = "BROADBAND SYSTEMS CORPORATION  LTD"
$ws.Range("E5").Value = "Remera, Gisimenti Airport Road (kn5 Rda), Opposite ChezLando"
$ws.Range("F5").Value = "Remera, Gisimenti Airport Road (kn5 Rda), Opposite ChezLando"
$ws.Range("G5").Value = "Remera, Gisimenti Airport Road (kn5 Rda), Opposite ChezLando"
$ws.Range("H5").Value = "7229 KIGALI, RWANDA"
$ws.Range("J5").Value = "gilbert.kayinamura@bsc.rw"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "4141"
$ws.Range("M5").Value = "www.bsc.rw"
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "1198080006094070"
$ws.Range("P5").Value = "BK"
$ws.Range("Q5").NumberFormat = "@"
$ws.Range("Q5").Value = "101982714"

# Row 6
$ws.Range("A6").Value = "3114072400000X"
$ws.Range("B6").Value = "CONNECT GLOBAL BUSINESS COMPANY Ltd"
$ws.Range("D6").Value = "CONNECT GLOBAL BUSINESS COMPANY Ltd"
$ws.Range("E6").Value = "KIGALI - NYARUGENGE"
$ws.Range("F6").Value = "KIGALI - NYARUGENGE"
$ws.Range("G6").Value = "KIGALI - NYARUGENGE"
$ws.Range("H6").Value = "POB:5564 Kigali"
$ws.Range("J6").Value = "niyafeos@yahoo.fr"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "788757320"
$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value = "1197980004514176"
$ws.Range("O6").Value = "RWANDA"
$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = "106907027"

# Row 7
$ws.Range("A7").Value = "3100785000000J"
$ws.Range("B7").Value = "INTERSEC SECURITY COMPANY  LTD ISCO"
$ws.Range("D7").Value = "INTERSEC SECURITY COMPANY  LTD ISCO"
$ws.Range("E7").Value = "Kanogo Gikondo Kicukiro"
$ws.Range("F7").Value = "Kanogo Gikondo Kicukiro"
$ws.Range("G7").Value = "Kanogo Gikondo Kicukiro"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "2146"
$ws.Range("J7").Value = "dnyangezi@isco.co.rw"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "0788310020"
$ws.Range("M7").Value = "www.isco.co.rw"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "1198770170398131"
$ws.Range("O7").Value = "RWANDA"
$ws.Range("P7").Value = "BK"
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "100076965"

# Row 8
$ws.Range("A8").Value = "3109596400000L"
$ws.Range("B8").Value = "NEPOMSCENE BUSINESS COMPANY  LTD"
$ws.Range("D8").Value = "NEPOMSCENE BUSINESS COMPANY  LTD"
$ws.Range("E8").Value = "NYAMAGABE District,Gasaka Sector"
$ws.Range("F8").Value = "NYAMAGABE District,Gasaka Sector"
$ws.Range("G8").Value = "NYAMAGABE District,Gasaka Sector"
$ws.Range("H8").ClearContents()
$ws.Range("J8").Value = "nbchvgmn@gmail.com"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "0788475217"
$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = "1198380172943184"
$ws.Range("P8").Value = "EQUITY BANK"
$ws.Range("Q8").NumberFormat = "@"
$ws.Range("Q8").Value = "103496546"

